$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025-04-03T19:32:30.695Z", "pravin", "pravin@gmail.com", "8363788798", "pravin"),
    @("2025-04-03T19:32:34.490Z", "pravin", "pravin@gmail.com", "8363788798", "pravin"),
    @("2025-04-03T19:32:37.592Z", "pravin", "pravin@gmail.com", "8363788798", "pravin"),
    @("2025-04-03T19:34:12.873Z", "pravin", "pravin@gmail.com", "8363788798", "fff")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($row, $c + 1)
        $value = $rowData[$c]
        if ($c -eq 3) {
            # Phone column looks numeric - force text storage like the source data,
            # then clear the quote-prefix style so no extra formatting is applied.
            $cell.Value = "'" + $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
